# Auto-generated: apply cell value corrections to the profit/pricing tables
# across all eight job sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) per the
# scheduled market-data refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 3960
$ws.Range("J40").Value = 4066.6667
$ws.Range("L40").Value = 4066.6667
$ws.Range("N40").Value = -4416.6667
$ws.Range("H43").Value = 1217.375
$ws.Range("I43").Value = 1267.8334
$ws.Range("K43").Value = 1267.8334
$ws.Range("M43").Value = -1198.8334
$ws.Range("H69").Value = 6000
$ws.Range("J69").Value = 6000
$ws.Range("L69").Value = 18000
$ws.Range("N69").Value = -19748
$ws.Range("H72").Value = 6000
$ws.Range("J72").Value = 6000
$ws.Range("L72").Value = 54000
$ws.Range("N72").Value = -62736
$ws.Range("H92").Value = 1697.5
$ws.Range("I92").Value = 2400
$ws.Range("K92").Value = 2400
$ws.Range("M92").Value = -1152
$ws.Range("H100").Value = 4606.8
$ws.Range("I100").Value = 1999.6666
$ws.Range("K100").Value = 1999.6666
$ws.Range("M100").Value = -1458.6666
$ws.Range("H138").Value = 6555.3335
$ws.Range("J138").Value = 5909.278
$ws.Range("L138").Value = 17727.834
$ws.Range("N138").Value = -28007.834

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 281.66666
$ws.Range("I4").Value = 281.66666
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 281.66666
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -165.66666
$ws.Range("N4").Value = $null
$ws.Range("H5").Value = 6299.5
$ws.Range("I5").Value = 3399.3333
$ws.Range("J5").Value = 15000
$ws.Range("K5").Value = 3399.3333
$ws.Range("L5").Value = 15000
$ws.Range("M5").Value = -3287.3333
$ws.Range("N5").Value = -15224
$ws.Range("H6").Value = 0
$ws.Range("I6").Value = 0
$ws.Range("K6").Value = 0
$ws.Range("M6").Value = $null
$ws.Range("H32").Value = 212000.98
$ws.Range("I32").Value = 212000.98
$ws.Range("K32").Value = 212000.98
$ws.Range("M32").Value = -211713.98
$ws.Range("H45").Value = 1751.8572
$ws.Range("I45").Value = 1614
$ws.Range("K45").Value = 1614
$ws.Range("M45").Value = -1237
$ws.Range("H101").Value = 391150.25
$ws.Range("J101").Value = 391150.25
$ws.Range("L101").Value = 391150.25
$ws.Range("N101").Value = -397640.25
$ws.Range("H132").Value = 1473026.9
$ws.Range("I132").Value = 1668857.1
$ws.Range("K132").Value = 5006571.300000001
$ws.Range("M132").Value = -5004041.300000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 6299.5
$ws.Range("I4").Value = 3399.3333
$ws.Range("J4").Value = 15000
$ws.Range("K4").Value = 3399.3333
$ws.Range("L4").Value = 15000
$ws.Range("M4").Value = -3284.3333
$ws.Range("N4").Value = -15230
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("M22").Value = $null
$ws.Range("H86").Value = 3813
$ws.Range("I86").Value = 3143
$ws.Range("K86").Value = 3143
$ws.Range("M86").Value = -2020
$ws.Range("H89").Value = 3813
$ws.Range("I89").Value = 3143
$ws.Range("K89").Value = 15715
$ws.Range("M89").Value = -10099
$ws.Range("H105").Value = 3954.125
$ws.Range("I105").Value = 1341.7778
$ws.Range("K105").Value = 1341.7778
$ws.Range("M105").Value = 405.2221999999999
$ws.Range("H134").Value = 15175719
$ws.Range("I134").Value = 28740.666
$ws.Range("K134").Value = 86221.99800000001
$ws.Range("M134").Value = -83686.99800000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 17037.729
$ws.Range("I7").Value = 45494
$ws.Range("J7").Value = 117.78378
$ws.Range("K7").Value = 45494
$ws.Range("L7").Value = 117.78378
$ws.Range("M7").Value = -45381
$ws.Range("N7").Value = -343.78378
$ws.Range("H22").Value = 74259.14
$ws.Range("I22").Value = 272.33334
$ws.Range("K22").Value = 272.33334
$ws.Range("M22").Value = 77.66665999999998
$ws.Range("H28").Value = 191881
$ws.Range("I28").Value = 300000
$ws.Range("K28").Value = 300000
$ws.Range("M28").Value = -299755
$ws.Range("H31").Value = 2648492.5
$ws.Range("I31").Value = 3270921.8
$ws.Range("J31").Value = 3168.25
$ws.Range("K31").Value = 3270921.8
$ws.Range("L31").Value = 3168.25
$ws.Range("M31").Value = -3270626.8
$ws.Range("N31").Value = -3758.25
$ws.Range("H34").Value = 2648492.5
$ws.Range("I34").Value = 3270921.8
$ws.Range("J34").Value = 3168.25
$ws.Range("K34").Value = 3270921.8
$ws.Range("L34").Value = 3168.25
$ws.Range("M34").Value = -3270719.8
$ws.Range("N34").Value = -3572.25
$ws.Range("H132").Value = 2478.3845
$ws.Range("I132").Value = 2170.0454
$ws.Range("J132").Value = 4174.25
$ws.Range("K132").Value = 6510.1362
$ws.Range("L132").Value = 12522.75
$ws.Range("M132").Value = -3980.1362
$ws.Range("N132").Value = -17582.75
$ws.Range("H134").Value = 2978.6943
$ws.Range("I134").Value = 2491.6875
$ws.Range("J134").Value = 6874.75
$ws.Range("K134").Value = 7475.0625
$ws.Range("L134").Value = 20624.25
$ws.Range("M134").Value = -4940.0625
$ws.Range("N134").Value = -25694.25
$ws.Range("H141").Value = 245138.47
$ws.Range("J141").Value = 253201.78
$ws.Range("L141").Value = 253201.78
$ws.Range("N141").Value = -263561.78

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H99").Value = 8040.5
$ws.Range("I99").Value = 2109.5715
$ws.Range("J99").Value = 13971.429
$ws.Range("K99").Value = 6328.7145
$ws.Range("L99").Value = 41914.287
$ws.Range("M99").Value = -4082.7145
$ws.Range("N99").Value = -46406.287
$ws.Range("H107").Value = 2401.8333
$ws.Range("J107").Value = 3012.3333
$ws.Range("L107").Value = 9036.999899999999
$ws.Range("N107").Value = -12876.9999
$ws.Range("H113").Value = 3608.3
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 3608.3
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 10824.9
$ws.Range("M113").Value = $null
$ws.Range("N113").Value = -15164.9
$ws.Range("H122").Value = 1153858.6
$ws.Range("J122").Value = 2015.75
$ws.Range("L122").Value = 18141.75
$ws.Range("N122").Value = -23041.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 5316087
$ws.Range("I2").Value = 5941444
$ws.Range("J2").Value = 549.5
$ws.Range("K2").Value = 5941444
$ws.Range("L2").Value = 549.5
$ws.Range("M2").Value = -5941331
$ws.Range("N2").Value = -775.5
$ws.Range("H34").Value = 69999
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 69999
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 69999
$ws.Range("M34").Value = $null
$ws.Range("N34").Value = -70535
$ws.Range("H76").Value = 69999
$ws.Range("I76").Value = 0
$ws.Range("J76").Value = 69999
$ws.Range("K76").Value = 0
$ws.Range("L76").Value = 69999
$ws.Range("M76").Value = $null
$ws.Range("N76").Value = -70629
$ws.Range("H79").Value = 69999
$ws.Range("I79").Value = 0
$ws.Range("J79").Value = 69999
$ws.Range("K79").Value = 0
$ws.Range("L79").Value = 69999
$ws.Range("M79").Value = $null
$ws.Range("N79").Value = -72183
$ws.Range("H132").Value = 14497.966
$ws.Range("I132").Value = 19480.54
$ws.Range("K132").Value = 58441.62
$ws.Range("M132").Value = -55911.62

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H42").Value = 44514
$ws.Range("J42").Value = 60028
$ws.Range("L42").Value = 60028
$ws.Range("N42").Value = -61154
$ws.Range("H46").Value = 5641.143
$ws.Range("J46").Value = 5844.3076
$ws.Range("L46").Value = 5844.3076
$ws.Range("N46").Value = -6220.3076
$ws.Range("H49").Value = 44514
$ws.Range("J49").Value = 60028
$ws.Range("L49").Value = 60028
$ws.Range("N49").Value = -60322
$ws.Range("H82").Value = 52216.9
$ws.Range("I82").Value = 64500
$ws.Range("J82").Value = 3084.5
$ws.Range("K82").Value = 64500
$ws.Range("L82").Value = 3084.5
$ws.Range("M82").Value = -64139
$ws.Range("N82").Value = -3806.5
$ws.Range("H85").Value = 52216.9
$ws.Range("I85").Value = 64500
$ws.Range("J85").Value = 3084.5
$ws.Range("K85").Value = 64500
$ws.Range("L85").Value = 3084.5
$ws.Range("M85").Value = -63252
$ws.Range("N85").Value = -5580.5
$ws.Range("H136").Value = 8936304
$ws.Range("I136").Value = 5686321
$ws.Range("K136").Value = 17058963
$ws.Range("M136").Value = -17056413

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 8697278
$ws.Range("I136").Value = 1978017.9
$ws.Range("J136").Value = 33334568
$ws.Range("K136").Value = 5934053.699999999
$ws.Range("L136").Value = 100003704
$ws.Range("M136").Value = -5931503.699999999
$ws.Range("N136").Value = -100008804
